$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83, shifting existing rows 83-128 down to 84-129
$ws.Rows.Item(83).Insert()

# Populate the new row 83 with the new weekly record
$ws.Range("A83").Value() = 11
$ws.Range("B83").Value() = "Vega Monumental Concepción"
$ws.Range("C83").Value() = "Bíobío"
$ws.Range("D83").Value() = 44813
$ws.Range("E83").Value() = 8
$ws.Range("F83").Value() = "Fruta"
$ws.Range("G83").Value() = 100108
$ws.Range("H83").Value() = "Tropicales y subtropicales"
$ws.Range("I83").Value() = 100108002
$ws.Range("J83").Value() = "Mango"
$ws.Range("K83").Value() = "Sin especificar"
$ws.Range("L83").Value() = "Primera"
$ws.Range("M83").Value() = 110
$ws.Range("N83").Value() = 9000
$ws.Range("O83").Value() = 9500
$ws.Range("P83").Value() = 9273
$ws.Range("Q83").Value() = "$/bandeja 4 kilos"
$ws.Range("R83").Value() = "Perú"
$ws.Range("S83").Value() = 2318
$ws.Range("T83").Value() = 4
